$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 2: B2 becomes EURUSD, clear C2 and D2
$ws.Range("B2").Value = "EURUSD"
$ws.Range("C2").ClearContents()
$ws.Range("D2").ClearContents()

# Update row 3: B3 becomes the new spot FX value, clear C3 and D3
$ws.Range("B3").Value = 1.1212599999999999
$ws.Range("C3").ClearContents()
$ws.Range("D3").ClearContents()

# Update row 4: B4 becomes the new vol curve name, clear C4 and D4
$ws.Range("B4").Value = "EURUSD VOL 8Y COTERMINAL 31122019"
$ws.Range("C4").ClearContents()
$ws.Range("D4").ClearContents()

# Update the selection on sheet1 to F10
[void]$ws.Range("F10").Select()
